$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.599.86'
$ws.Range('E2').Value = '  -2.22%  '

$ws.Range('D3').Value = '2.959.92'
$ws.Range('E3').Value = '  -3.30%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '496.97'
$ws.Range('E5').Value = '  -5.62%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.80'
$ws.Range('E6').Value = '  -5.92%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.425'
$ws.Range('E8').Value = '  -5.08%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.18'
$ws.Range('E9').Value = '  -5.93%  '

$ws.Range('E10').Value = '  -6.21%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.353'
$ws.Range('E11').Value = '  -4.77%  '

$ws.Range('D12').Value = '3.462.00'
$ws.Range('E12').Value = '  -3.54%  '

$ws.Range('E13').Value = '  -3.25%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.94'
$ws.Range('E14').Value = '  -5.04%  '

$ws.Range('E15').Value = '  -8.07%  '

$ws.Range('D16').Value = '56.546.49'
$ws.Range('E16').Value = '  -2.30%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.01'
$ws.Range('E17').Value = '  -3.26%  '

$ws.Range('D18').Value = '2.953.04'
$ws.Range('E18').Value = '  -3.58%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.50'
$ws.Range('E19').Value = '  -5.13%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.77'
$ws.Range('E20').Value = '  -5.39%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '317.17'
$ws.Range('E21').Value = '  -7.08%  '

$ws.Range('E22').Value = '  +0.10%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.73'
$ws.Range('E23').Value = '  +0.90%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.486'
$ws.Range('E24').Value = '  -3.46%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '62.48'
$ws.Range('E25').Value = '  -3.82%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.20%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.162'
$ws.Range('E27').Value = '  -4.66%  '

$ws.Range('D28').Value = '0.0₃0871'
$ws.Range('E28').Value = '  -10.99%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.49'
$ws.Range('E29').Value = '  -7.03%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.00'
$ws.Range('E30').Value = '  -5.48%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.75'
$ws.Range('E31').Value = '  -5.99%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.93'
$ws.Range('E32').Value = '  -5.42%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.14'
$ws.Range('E33').Value = '  -8.48%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '152.99'
$ws.Range('E34').Value = '  -2.69%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.49'
$ws.Range('E35').Value = '  -6.39%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.70'
$ws.Range('E36').Value = '  -4.71%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.21'
$ws.Range('E37').Value = '  -9.11%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.85'
$ws.Range('E38').Value = '  -9.31%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0656'
$ws.Range('E39').Value = '  -7.05%  '

$ws.Range('D40').Value = '2.987.72'
$ws.Range('E40').Value = '  -3.56%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '37.34'
$ws.Range('E41').Value = '  -1.43%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.20%  '

$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.69'
$ws.Range('E43').Value = '  -6.43%  '

$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.640'
$ws.Range('E44').Value = '  -4.00%  '

$ws.Range('D45').Value = '2.150.95'
$ws.Range('E45').Value = '  -8.14%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.36'
$ws.Range('E46').Value = '  -8.41%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.86'
$ws.Range('E47').Value = '  -3.12%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.921'
$ws.Range('E48').Value = '  -10.89%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0232'
$ws.Range('E49').Value = '  -5.39%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.04'
$ws.Range('E50').Value = '  -5.75%  '

$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.74'
$ws.Range('E51').Value = '  -13.43%  '
